$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = 3
$ws.Range("D5").Value = 89
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 13
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 13
$ws.Range("D13").Value = 3
$ws.Range("D14").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("D20").Value = 89
$ws.Range("D21").Value = 13
$ws.Range("D22").Value = 13
$ws.Range("D24").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("D26").Value = 21
$ws.Range("D27").Value = 8
$ws.Range("D28").Value = 8
$ws.Range("D29").Value = 8
$ws.Range("D30").Value = 8
$ws.Range("D31").Value = 8
$ws.Range("D32").Value = 8
$ws.Range("D33").Value = 8
$ws.Range("D34").Value = 8
$ws.Range("D35").Value = 8
$ws.Range("D36").Value = 8
$ws.Range("D37").Value = 8
$ws.Range("D38").Value = 8
$ws.Range("D39").Value = 8
$ws.Range("D40").Value = 8
$ws.Range("D41").Value = 8
$ws.Range("D42").Value = 8
$ws.Range("D43").Value = 8
$ws.Range("D44").Value = 8
$ws.Range("D45").Value = 8
$ws.Range("D46").Value = 8
$ws.Range("D47").Value = 8
$ws.Range("D48").Value = 8
$ws.Range("D49").Value = 8
$ws.Range("D50").Value = 8
$ws.Range("D51").Value = 8
$ws.Range("D52").Value = 8
$ws.Range("D53").Value = 8
$ws.Range("D54").Value = 8
$ws.Range("D55").Value = 8
$ws.Range("D56").Value = 8
$ws.Range("D57").Value = 8
$ws.Range("D58").Value = 8
$ws.Range("D59").Value = 8
$ws.Range("D60").Value = 8
$ws.Range("D61").Value = 8
$ws.Range("D62").Value = 3
$ws.Range("D63").Value = 3
$ws.Range("D64").Value = 13
$ws.Range("D65").Value = 5
$ws.Range("D66").Value = 5
$ws.Range("D67").Value = 5
$ws.Range("D68").Value = 21
$ws.Range("D69").Value = 21
$ws.Range("D70").Value = 3
$ws.Range("D71").Value = 5
$ws.Range("D72").Value = 8
$ws.Range("D73").Value = 5
$ws.Range("D74").Value = 5
$ws.Range("D75").Value = 5
$ws.Range("D77").Value = 8
$ws.Range("D78").Value = 89
$ws.Range("D79").Value = 3
$ws.Range("D80").Value = 5
$ws.Range("D81").Value = 3
$ws.Range("D82").Value = 34
$ws.Range("D83").Value = 21
$ws.Range("D84").Value = 8
$ws.Range("D85").Value = 3
$ws.Range("D86").Value = 5
$ws.Range("D87").Value = 89
$ws.Range("D88").Value = 21
$ws.Range("D89").Value = 3
$ws.Range("D90").Value = 5
$ws.Range("D91").Value = 55
$ws.Range("D92").Value = 8
$ws.Range("D93").Value = 8
$ws.Range("D94").Value = 21
$ws.Range("D95").Value = 5
$ws.Range("D96").Value = 89
$ws.Range("D97").Value = 3
$ws.Range("D98").Value = 3
$ws.Range("D99").Value = 5
$ws.Range("D100").Value = 5
$ws.Range("D101").Value = 5
$ws.Range("D102").Value = 3
$ws.Range("D103").Value = 21
$ws.Range("D104").Value = 89
$ws.Range("D105").Value = 89
$ws.Range("D106").Value = 21
$ws.Range("D107").Value = 3
$ws.Range("D108").Value = 3
$ws.Range("D109").Value = 5
$ws.Range("D110").Value = 8
$ws.Range("D111").Value = 8
$ws.Range("D112").Value = 3
$ws.Range("D113").Value = 5
$ws.Range("D114").Value = 8
$ws.Range("D115").Value = 5
$ws.Range("D116").Value = 21
$ws.Range("D117").Value = 21
$ws.Range("D118").Value = 21
$ws.Range("D119").Value = 8
$ws.Range("D120").Value = 3
$ws.Range("D121").Value = 3
$ws.Range("D122").Value = 3
$ws.Range("D123").Value = 5
$ws.Range("D124").Value = 8
$ws.Range("D125").Value = 8
$ws.Range("D126").Value = 8
$ws.Range("D127").Value = 89
$ws.Range("D128").Value = 3
$ws.Range("D129").Value = 3
$ws.Range("D130").Value = 34
$ws.Range("D131").Value = 8
$ws.Range("D132").Value = 3
$ws.Range("D133").Value = 3
$ws.Range("D134").Value = 5
$ws.Range("D135").Value = 3
$ws.Range("D136").Value = 21
$ws.Range("D137").Value = 89
$ws.Range("D138").Value = 3
$ws.Range("D139").Value = 8
$ws.Range("D140").Value = 13
$ws.Range("D141").Value = 3
$ws.Range("D142").Value = 21
$ws.Range("D143").Value = 21
$ws.Range("D144").Value = 21
$ws.Range("D145").Value = 21
$ws.Range("D146").Value = 21
$ws.Range("D147").Value = 21
$ws.Range("D148").Value = 21
$ws.Range("D149").Value = 21
$ws.Range("D150").Value = 21
$ws.Range("D151").Value = 21
$ws.Range("D152").Value = 21
$ws.Range("D153").Value = 21
$ws.Range("D154").Value = 21
$ws.Range("D155").Value = 21
$ws.Range("D156").Value = 21
$ws.Range("D157").Value = 21
$ws.Range("D158").Value = 21
$ws.Range("D159").Value = 21
$ws.Range("D160").Value = 21
$ws.Range("D161").Value = 21
$ws.Range("D162").Value = 21
$ws.Range("D163").Value = 21
$ws.Range("D164").Value = 21
$ws.Range("D165").Value = 21
$ws.Range("D166").Value = 21
$ws.Range("D167").Value = 8
$ws.Range("D168").Value = 3
$ws.Range("D169").Value = 3
$ws.Range("D170").Value = 3
$ws.Range("D171").Value = 21
$ws.Range("D172").Value = 5
$ws.Range("D173").Value = 5
$ws.Range("D174").Value = 89
$ws.Range("D175").Value = 5
$ws.Range("D176").Value = 34
$ws.Range("D177").Value = 3
$ws.Range("D178").Value = 3
$ws.Range("D179").Value = 55
$ws.Range("D180").Value = 55
$ws.Range("D181").Value = 55
$ws.Range("D182").Value = 55
$ws.Range("D183").Value = 55
$ws.Range("D184").Value = 55
$ws.Range("D185").Value = 55
$ws.Range("D186").Value = 55
$ws.Range("D187").Value = 55
$ws.Range("D188").Value = 55
$ws.Range("D189").Value = 55
$ws.Range("D190").Value = 55
$ws.Range("D191").Value = 55
$ws.Range("D192").Value = 55
$ws.Range("D193").Value = 55
$ws.Range("D194").Value = 3
$ws.Range("D195").Value = 8
$ws.Range("D196").Value = 8
$ws.Range("D197").Value = 55
$ws.Range("D198").Value = 34
$ws.Range("D199").Value = 3
$ws.Range("D200").Value = 3
$ws.Range("D201").Value = 3
$ws.Range("D202").Value = 13
$ws.Range("D203").Value = 3
$ws.Range("D204").Value = 5
$ws.Range("D205").Value = 89
$ws.Range("D206").Value = 3
$ws.Range("D207").Value = 3
$ws.Range("D208").Value = 3
$ws.Range("D209").Value = 3
$ws.Range("D210").Value = 3
$ws.Range("D211").Value = 3
$ws.Range("D212").Value = 5
$ws.Range("D213").Value = 3
$ws.Range("D214").Value = 5
$ws.Range("D215").Value = 5
$ws.Range("D216").Value = 5
$ws.Range("D217").Value = 3
$ws.Range("D218").Value = 5
$ws.Range("D219").Value = 5
$ws.Range("D220").Value = 5
$ws.Range("D221").Value = 5
$ws.Range("D222").Value = 5
$ws.Range("D223").Value = 5
$ws.Range("D224").Value = 5
$ws.Range("D225").Value = 5
$ws.Range("D226").Value = 5
$ws.Range("D227").Value = 5
$ws.Range("D228").Value = 5
$ws.Range("D229").Value = 21
$ws.Range("D230").Value = 34
$ws.Range("D231").Value = 34
$ws.Range("D232").Value = 8
$ws.Range("D233").Value = 13
$ws.Range("D234").Value = 5
$ws.Range("D235").Value = 3
$ws.Range("D236").Value = 3
$ws.Range("D237").Value = 8
$ws.Range("D238").Value = 21
$ws.Range("D239").Value = 55
$ws.Range("D240").Value = 13
$ws.Range("D241").Value = 3
$ws.Range("D242").Value = 3
$ws.Range("D243").Value = 3
$ws.Range("D244").Value = 3
$ws.Range("D245").Value = 3
$ws.Range("D246").Value = 8
$ws.Range("D247").Value = 8
$ws.Range("D248").Value = 89
$ws.Range("D249").Value = 89
$ws.Range("D250").Value = 89
$ws.Range("D251").Value = 89
$ws.Range("D252").Value = 89
$ws.Range("D253").Value = 89
$ws.Range("D254").Value = 55
$ws.Range("D255").Value = 3
$ws.Range("D256").Value = 13
$ws.Range("D257").Value = 3
$ws.Range("D258").Value = 3
$ws.Range("D259").Value = 13
$ws.Range("D260").Value = 3
$ws.Range("D261").Value = 13
$ws.Range("D262").Value = 5
$ws.Range("D263").Value = 8
$ws.Range("D264").Value = 21
$ws.Range("D265").Value = 3
$ws.Range("D266").Value = 8
$ws.Range("D267").Value = 8
$ws.Range("D268").Value = 8
$ws.Range("D269").Value = 55
$ws.Range("D270").Value = 3
$ws.Range("D271").Value = 5
$ws.Range("D272").Value = 8
$ws.Range("D273").Value = 5
$ws.Range("D274").Value = 8
$ws.Range("D275").Value = 5
$ws.Range("D276").Value = 3
$ws.Range("D277").Value = 8
$ws.Range("D278").Value = 5
$ws.Range("D279").Value = 21
$ws.Range("D280").Value = 5
$ws.Range("D281").Value = 3
$ws.Range("D282").Value = 8
$ws.Range("D283").Value = 13
$ws.Range("D285").Value = 3
$ws.Range("D286").Value = 3
$ws.Range("D287").Value = 5
$ws.Range("D288").Value = 3
$ws.Range("D289").Value = 13
$ws.Range("D290").Value = 13
$ws.Range("D291").Value = 13
$ws.Range("D292").Value = 13
$ws.Range("D293").Value = 13
$ws.Range("D294").Value = 13
$ws.Range("D295").Value = 13
$ws.Range("D296").Value = 13
$ws.Range("D297").Value = 13
$ws.Range("D298").Value = 13
$ws.Range("D299").Value = 13
$ws.Range("D300").Value = 13
$ws.Range("D301").Value = 13
$ws.Range("D302").Value = 13
$ws.Range("D303").Value = 13
$ws.Range("D304").Value = 13
$ws.Range("D305").Value = 13
$ws.Range("D306").Value = 13
$ws.Range("D307").Value = 13
$ws.Range("D308").Value = 13
$ws.Range("D309").Value = 13
$ws.Range("D310").Value = 13
$ws.Range("D311").Value = 13
$ws.Range("D312").Value = 13
$ws.Range("D313").Value = 13
$ws.Range("D314").Value = 13
$ws.Range("D315").Value = 13
$ws.Range("D316").Value = 13
$ws.Range("D317").Value = 13
$ws.Range("D318").Value = 13
$ws.Range("D319").Value = 13
$ws.Range("D320").Value = 13
$ws.Range("D321").Value = 13
$ws.Range("D322").Value = 13
$ws.Range("D323").Value = 13
$ws.Range("D324").Value = 13
$ws.Range("D325").Value = 13
$ws.Range("D326").Value = 13
$ws.Range("D327").Value = 13
$ws.Range("D328").Value = 13
$ws.Range("D329").Value = 13
$ws.Range("D330").Value = 13
$ws.Range("D331").Value = 13
$ws.Range("D332").Value = 13
$ws.Range("D333").Value = 8
$ws.Range("D334").Value = 89
$ws.Range("D335").Value = 3
$ws.Range("D336").Value = 3
$ws.Range("D337").Value = 13
$ws.Range("D338").Value = 3
$ws.Range("D339").Value = 8
$ws.Range("D340").Value = 89
$ws.Range("D341").Value = 5
$ws.Range("D342").Value = 34
$ws.Range("D343").Value = 13
$ws.Range("D344").Value = 5
$ws.Range("D345").Value = 3
$ws.Range("D346").Value = 8
$ws.Range("D347").Value = 13
$ws.Range("D348").Value = 13
$ws.Range("D349").Value = 13
$ws.Range("D350").Value = 8
$ws.Range("D351").Value = 21
$ws.Range("D352").Value = 89
$ws.Range("D353").Value = 89
$ws.Range("D354").Value = 21
$ws.Range("D355").Value = 8
$ws.Range("D356").Value = 89
$ws.Range("D357").Value = 89
$ws.Range("D358").Value = 13
$ws.Range("D359").Value = 5
$ws.Range("D360").Value = 8
$ws.Range("D361").Value = 8
$ws.Range("D362").Value = 89
$ws.Range("D363").Value = 5
$ws.Range("D364").Value = 3
$ws.Range("D365").Value = 3
$ws.Range("D366").Value = 3
$ws.Range("D367").Value = 5
$ws.Range("D368").Value = 3
$ws.Range("D369").Value = 8
$ws.Range("D370").Value = 8
$ws.Range("D371").Value = 8
$ws.Range("D372").Value = 3
$ws.Range("D373").Value = 3
$ws.Range("D374").Value = 5
$ws.Range("D375").Value = 5
$ws.Range("D376").Value = 5
$ws.Range("D377").Value = 34
$ws.Range("D378").Value = 3
$ws.Range("D379").Value = 5
$ws.Range("D380").Value = 5
$ws.Range("D381").Value = 5
$ws.Range("D382").Value = 5
$ws.Range("D383").Value = 5
$ws.Range("D384").Value = 3
$ws.Range("D385").Value = 3
$ws.Range("D386").Value = 8
$ws.Range("D387").Value = 3
$ws.Range("D388").Value = 21
$ws.Range("D389").Value = 3
$ws.Range("D390").Value = 3
$ws.Range("D391").Value = 3
$ws.Range("D392").Value = 13
$ws.Range("D393").Value = 21
$ws.Range("D394").Value = 5
$ws.Range("D395").Value = 13
$ws.Range("D396").Value = 3
$ws.Range("D397").Value = 3
$ws.Range("D398").Value = 3
$ws.Range("D399").Value = 3
$ws.Range("D400").Value = 3
$ws.Range("D401").Value = 3
$ws.Range("D402").Value = 3
$ws.Range("D403").Value = 3
$ws.Range("D404").Value = 3
$ws.Range("D405").Value = 3
$ws.Range("D406").Value = 3
$ws.Range("D407").Value = 3
$ws.Range("D408").Value = 3
$ws.Range("D409").Value = 3
$ws.Range("D410").Value = 3
$ws.Range("D411").Value = 3
$ws.Range("D412").Value = 3
$ws.Range("D413").Value = 3
$ws.Range("D414").Value = 3
$ws.Range("D415").Value = 3
$ws.Range("D416").Value = 3
$ws.Range("D417").Value = 3
$ws.Range("D418").Value = 3
$ws.Range("D419").Value = 3
$ws.Range("D420").Value = 3
$ws.Range("D421").Value = 3
$ws.Range("D422").Value = 3
$ws.Range("D423").Value = 3
$ws.Range("D424").Value = 3
$ws.Range("D425").Value = 3
$ws.Range("D426").Value = 3
$ws.Range("D427").Value = 3
$ws.Range("D428").Value = 3
$ws.Range("D429").Value = 3
$ws.Range("D430").Value = 3
$ws.Range("D431").Value = 3
$ws.Range("D432").Value = 3
$ws.Range("D433").Value = 3
$ws.Range("D434").Value = 3
$ws.Range("D435").Value = 3
$ws.Range("D436").Value = 3
$ws.Range("D437").Value = 3
$ws.Range("D438").Value = 3
$ws.Range("D439").Value = 3
$ws.Range("D440").Value = 3
$ws.Range("D441").Value = 3
$ws.Range("D442").Value = 3
$ws.Range("D443").Value = 3
$ws.Range("D444").Value = 3
$ws.Range("D445").Value = 3
$ws.Range("D446").Value = 3
$ws.Range("D447").Value = 3
$ws.Range("D448").Value = 3
$ws.Range("D449").Value = 3
$ws.Range("D450").Value = 3
$ws.Range("D451").Value = 3
$ws.Range("D452").Value = 3
$ws.Range("D453").Value = 3
$ws.Range("D454").Value = 3
$ws.Range("D455").Value = 3
$ws.Range("D456").Value = 3
$ws.Range("D457").Value = 3
$ws.Range("D458").Value = 3
$ws.Range("D459").Value = 3
$ws.Range("D460").Value = 3
$ws.Range("D461").Value = 3
$ws.Range("D462").Value = 3
$ws.Range("D463").Value = 3
$ws.Range("D464").Value = 3
$ws.Range("D465").Value = 3
$ws.Range("D466").Value = 3
$ws.Range("D467").Value = 3
$ws.Range("D468").Value = 3
$ws.Range("D469").Value = 3
$ws.Range("D470").Value = 3
$ws.Range("D471").Value = 5
$ws.Range("D472").Value = 5
$ws.Range("D473").Value = 5
$ws.Range("D474").Value = 3
$ws.Range("D475").Value = 5
$ws.Range("D476").Value = 13
$ws.Range("D477").Value = 3
$ws.Range("D478").Value = 34
$ws.Range("D479").Value = 8
$ws.Range("D480").Value = 3
$ws.Range("D481").Value = 5
$ws.Range("D482").Value = 3
$ws.Range("D483").Value = 5
$ws.Range("D484").Value = 3
$ws.Range("D485").Value = 3
$ws.Range("D486").Value = 3
$ws.Range("D487").Value = 8
$ws.Range("D488").Value = 21
$ws.Range("D489").Value = 13
